$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 10 ("Objetivos:" / "Objectives:") - the Portuguese body text is replaced
# by the professor name that used to live further down the sheet.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Range("C10").Value = "1304060 - Maria das Graças de Almeida Felipe"

# ---------------------------------------------------------------------------
# Remove the old rows 13-25 entirely, then rebuild rows 13-23 from scratch
# with the new content / layout.
# ---------------------------------------------------------------------------
$ws.Rows("13:25").Delete()

# Row 13 - "Programa resumido:" row (value reused from the "Ativação:" cell)
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Rows(13).RowHeight = 60

# Row 14 - "Short syllabus:" row
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Development of the course conclusion work under orientation of a leader professor, which must deal with subjects related to Biochemical Engineering."
$ws.Range("C14").Value = "Development of the course conclusion work under orientation of a leader professor, which must deal with subjects related to Biochemical Engineering."
$ws.Rows(14).RowHeight = 60

# Row 15 - "Programa:" row (value reused from the professor name cell above)
$ws.Range("A15").Value = "Programa:"
$ws.Range("B10").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C10").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Rows(15).RowHeight = 120

# Row 16 - "Syllabus:" row
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Elaboration of a monograph of course conclusion presenting: (1) the subject and its importance, (2) the objectives, (3) the bibliographic revision, (4) the scientific methodology, (5) the development of the project, (6) the analysis and discussion of the results, (7) the conclusion and recommendations for the future works and (8) references. The document must attend to the Brazilian Association of Technical Norms."
$ws.Range("C16").Value = "Elaboration of a monograph of course conclusion presenting: (1) the subject and its importance, (2) the objectives, (3) the bibliographic revision, (4) the scientific methodology, (5) the development of the project, (6) the analysis and discussion of the results, (7) the conclusion and recommendations for the future works and (8) references. The document must attend to the Brazilian Association of Technical Norms."
$ws.Rows(16).RowHeight = 120

# Row 17 - "Avaliação:" (label only, no content cells)
$ws.Range("A17").Value = "Avaliação:"

# Row 18 - "Método:" row
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "8853480 - Tatiane da Franca Silva"
$ws.Range("C18").Value = "8853480 - Tatiane da Franca Silva"
$ws.Rows(18).RowHeight = 60

# Row 19 - "Critério:" row
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Análise da monografia e defesa do trabalho perante uma banca de 3 examinadores, conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica."
$ws.Range("C19").Value = "Análise da monografia e defesa do trabalho perante uma banca de 3 examinadores, conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica."
$ws.Rows(19).RowHeight = 60

# Row 20 - "Norma de recuperação:" row
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota da disciplina será decidida pelos docentes da banca"
$ws.Range("C20").Value = "A nota da disciplina será decidida pelos docentes da banca"
$ws.Rows(20).RowHeight = 60

# Row 21 - "Bibliografia:" row
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Reapresentação do trabalho modificado para nova avaliação."
$ws.Range("C21").Value = "Reapresentação do trabalho modificado para nova avaliação."
$ws.Rows(21).RowHeight = 120

# Row 22 - "Requisitos:" (label only, no content cells)
$ws.Range("A22").Value = "Requisitos:"

# Row 23 - prerequisite course reference (no label cell)
$ws.Range("B23").Value = "LOT2056 -  Trabalho de Conclusão de Curso I  (Requisito)
"
$ws.Range("C23").Value = "LOT2056 -  Trabalho de Conclusão de Curso I  (Requisito)
"
$ws.Rows(23).RowHeight = 30

# ---------------------------------------------------------------------------
# Re-apply the column styles (A=bold/top, B=wrap/top, C=red/wrap/top) that the
# freshly written rows should carry, by copying formats from row 3, which
# keeps the original A/B/C column styling untouched by this edit.
# ---------------------------------------------------------------------------
$ws.Range("A3:C3").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A18:C18").PasteSpecial(-4122)
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Range("A21:C21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A3").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B3:C3").Copy()
$ws.Range("B23:C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
